# Auto update Excel log
$wb = $excel.ActiveWorkbook

function Set-RowData($ws, $row, $date, $time, $hour, $loc, $val, $status) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $date
    $cellA.Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $time
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $loc

    $cellE = $ws.Cells.Item($row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $val
    $cellE.Style = "Normal"

    $ws.Cells.Item($row, 6).Value = $status
}

$pir = $wb.Worksheets.Item("PIR")

Set-RowData $pir 305 "2026-01-30" "17:32:26" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 306 "2026-01-30" "17:32:27" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 307 "2026-01-30" "17:32:30" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 308 "2026-01-30" "17:32:35" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 309 "2026-01-30" "17:32:40" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 310 "2026-01-30" "17:32:45" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 311 "2026-01-30" "17:32:50" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 312 "2026-01-30" "17:32:55" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 313 "2026-01-30" "17:33:00" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 314 "2026-01-30" "17:33:05" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 315 "2026-01-30" "17:33:10" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 316 "2026-01-30" "17:33:15" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 317 "2026-01-30" "17:33:20" "17:00" "Bathroom" "No Motion" "Inactive"
Set-RowData $pir 318 "2026-01-30" "17:33:25" "17:00" "Bathroom" "No Motion" "Inactive"

$hum = $wb.Worksheets.Item("Humidity")

Set-RowData $hum 212 "2026-01-30" "17:32:26" "17:00" "Bathroom" "86.1%" "Active"
Set-RowData $hum 213 "2026-01-30" "17:32:36" "17:00" "Bathroom" "87.1%" "Active"
Set-RowData $hum 214 "2026-01-30" "17:32:41" "17:00" "Bathroom" "86.1%" "Active"
Set-RowData $hum 215 "2026-01-30" "17:32:56" "17:00" "Bathroom" "87.0%" "Active"
Set-RowData $hum 216 "2026-01-30" "17:33:01" "17:00" "Bathroom" "87.1%" "Active"
Set-RowData $hum 217 "2026-01-30" "17:33:11" "17:00" "Bathroom" "87.1%" "Active"
Set-RowData $hum 218 "2026-01-30" "17:33:16" "17:00" "Bathroom" "87.1%" "Active"
Set-RowData $hum 219 "2026-01-30" "17:33:21" "17:00" "Bathroom" "86.2%" "Active"
